# Suzie's IRA: update 2016 dividend figures (the 2017 sheet pulls from
# the 2016 sheet, and the Total/Prev Year sheet pulls from 2017, so those
# cascade automatically via formula recalculation).
$wb = $excel.ActiveWorkbook
$ws2016 = $wb.Worksheets.Item("2016")
$ws2016.Activate()

# Row 3 (EMR): price per share and Dec dividend payment updated
$ws2016.Range("D3").Value = 14.372999999999999
$ws2016.Range("S3").Value = 6.84

# Row 5 (KO): price per share and Dec dividend payment updated
$ws2016.Range("D5").Value = 40.033000000000001
$ws2016.Range("S5").Value = 13.9

# Row 8 (RDS.A): price per share and Dec dividend payment updated
$ws2016.Range("D8").Value = 18.715
$ws2016.Range("S8").Value = 17.260000000000002

# Match the author's last active cell selection on the 2016 sheet
$ws2016.Range("I14").Select()
